$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F4").Value = -7
$ws.Range("F5").Value = -6
$ws.Range("F6").Value = -2
$ws.Range("F8").Value = -11
$ws.Range("F9").Value = -7
$ws.Range("F11").Value = -3
